# Bugfix the naive forecaster component module:
# the YoY component forecast sheet drops its two oldest trailing
# history rows (23:24) and its newest trailing forecast-vintage
# column (BA), and recomputes the leading forecast values for
# every remaining forecast-vintage row / the date header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the two obsolete trailing rows (delete bottom-up so row
# numbers above stay stable) and the obsolete trailing column.
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(23).Delete()
$ws.Columns.Item(53).Delete()

# Recompute the forecast values that changed.
$arr = New-Object "object[,]" 1,51
$arr[0,0] = 39583
$arr[0,1] = 39765
$arr[0,2] = 39948
$arr[0,3] = 40130
$arr[0,4] = 40310
$arr[0,5] = 40494
$arr[0,6] = 40676
$arr[0,7] = 40862
$arr[0,8] = 41044
$arr[0,9] = 41228
$arr[0,10] = 41409
$arr[0,11] = 41592
$arr[0,12] = 41774
$arr[0,13] = 41957
$arr[0,14] = 42137
$arr[0,15] = 42321
$arr[0,16] = 42503
$arr[0,17] = 42689
$arr[0,18] = 42867
$arr[0,19] = 43053
$arr[0,20] = 43145
$arr[0,21] = 43235
$arr[0,22] = 43326
$arr[0,23] = 43418
$arr[0,24] = 43510
$arr[0,25] = 43600
$arr[0,26] = 43691
$arr[0,27] = 43783
$arr[0,28] = 43875
$arr[0,29] = 43966
$arr[0,30] = 44068
$arr[0,31] = 44159
$arr[0,32] = 44251
$arr[0,33] = 44341
$arr[0,34] = 44432
$arr[0,35] = 44525
$arr[0,36] = 44617
$arr[0,37] = 44706
$arr[0,38] = 44798
$arr[0,39] = 44890
$arr[0,40] = 44981
$arr[0,41] = 45071
$arr[0,42] = 45163
$arr[0,43] = 45254
$arr[0,44] = 45345
$arr[0,45] = 45436
$arr[0,46] = 45534
$arr[0,47] = 45618
$arr[0,48] = 45713
$arr[0,49] = 45800
$arr[0,50] = 45891
$ws.Range("B1:AZ1").Value2 = $arr

$arr = New-Object "object[,]" 1,3
$arr[0,0] = 1.782259294303912
$arr[0,1] = -0.08289353495386509
$arr[0,2] = -0.9921462019007898
$ws.Range("B3:D3").Value2 = $arr

$arr = New-Object "object[,]" 1,5
$arr[0,0] = 1.840084080815463
$arr[0,1] = 0.2227880631417101
$arr[0,2] = -0.3230348957779294
$arr[0,3] = -0.363786394693788
$arr[0,4] = 0.2467309912830284
$ws.Range("B4:F4").Value2 = $arr

$arr = New-Object "object[,]" 1,6
$arr[0,0] = $null
$arr[0,1] = -0.3323859540900087
$arr[0,2] = -0.3968757761298791
$arr[0,3] = -0.06959526544320083
$arr[0,4] = -0.2139598932957232
$arr[0,5] = 1.160201558804674
$ws.Range("C5:H5").Value2 = $arr

$arr = New-Object "object[,]" 1,6
$arr[0,0] = $null
$arr[0,1] = -0.1383094607783963
$arr[0,2] = -0.4829805246118979
$arr[0,3] = 0.2932139896134167
$arr[0,4] = 0.903223459378788
$arr[0,5] = 1.048604932640185
$ws.Range("E6:J6").Value2 = $arr

$arr = New-Object "object[,]" 1,6
$arr[0,0] = $null
$arr[0,1] = 0.2328086450296141
$arr[0,2] = 0.723280697834694
$arr[0,3] = 1.078804187516891
$arr[0,4] = 1.31837503023402
$arr[0,5] = 1.06837811337479
$ws.Range("G7:L7").Value2 = $arr

$arr = New-Object "object[,]" 1,6
$arr[0,0] = $null
$arr[0,1] = 1.043680860183693
$arr[0,2] = 1.185259818806217
$arr[0,3] = 1.25598608434605
$arr[0,4] = 1.658305347589661
$arr[0,5] = 1.435208340819005
$ws.Range("I8:N8").Value2 = $arr

$arr = New-Object "object[,]" 1,6
$arr[0,0] = $null
$arr[0,1] = 1.286157117685827
$arr[0,2] = 1.366951734963395
$arr[0,3] = 1.407107513712802
$arr[0,4] = 1.488472133572305
$arr[0,5] = 1.577589817310243
$ws.Range("K9:P9").Value2 = $arr

$arr = New-Object "object[,]" 1,6
$arr[0,0] = $null
$arr[0,1] = 1.366951734963395
$arr[0,2] = 1.441237482553381
$arr[0,3] = 1.464859320654099
$arr[0,4] = 1.644157643645183
$arr[0,5] = 1.979569114089963
$ws.Range("M10:R10").Value2 = $arr

$arr = New-Object "object[,]" 1,6
$arr[0,0] = $null
$arr[0,1] = 1.449109184169628
$arr[0,2] = 1.647704289169027
$arr[0,3] = 1.639881111696151
$arr[0,4] = 1.586470485311331
$arr[0,5] = 1.815212363528707
$ws.Range("O11:T11").Value2 = $arr

$arr = New-Object "object[,]" 1,8
$arr[0,0] = $null
$arr[0,1] = 1.603672482622964
$arr[0,2] = 1.652736919119047
$arr[0,3] = 1.806931013599544
$arr[0,4] = 1.974604558490256
$arr[0,5] = 2.140635848901895
$arr[0,6] = 2.284026378382942
$arr[0,7] = 2.336516087993035
$ws.Range("Q12:X12").Value2 = $arr

$arr = New-Object "object[,]" 1,11
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = 1.796783738056584
$arr[0,3] = 1.893938831797337
$arr[0,4] = 1.99288634244883
$arr[0,5] = 2.119133965447961
$arr[0,6] = 2.162438527487853
$arr[0,7] = 1.93172124148786
$arr[0,8] = 1.404530461900833
$arr[0,9] = 1.264761787657309
$arr[0,10] = 1.17909021197069
$ws.Range("R13:AB13").Value2 = $arr

$arr = New-Object "object[,]" 1,13
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = 2.083122398735981
$arr[0,4] = 2.071818698227212
$arr[0,5] = 1.970775328194052
$arr[0,6] = 1.887821778955101
$arr[0,7] = 1.798687504247187
$arr[0,8] = 1.636329093826605
$arr[0,9] = 1.008270799755984
$arr[0,10] = 0.5345697479163913
$arr[0,11] = 0.824608016336259
$arr[0,12] = -1.788000783651811
$ws.Range("T14:AF14").Value2 = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = 1.839044560557102
$arr[0,6] = 1.736580578290958
$arr[0,7] = 1.468343000448269
$arr[0,8] = 1.247274949485733
$arr[0,9] = 1.395219579261608
$arr[0,10] = -0.02261741485058977
$arr[0,11] = -1.119700950349478
$arr[0,12] = -2.680286313062752
$arr[0,13] = -1.4191429117966
$arr[0,14] = -1.098964423305859
$ws.Range("V15:AJ15").Value2 = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = 1.306001555547232
$arr[0,7] = -0.07129657200888317
$arr[0,8] = -0.6741130548587049
$arr[0,9] = -2.013357217277445
$arr[0,10] = -0.1125839228000469
$arr[0,11] = 1.055324027461602
$arr[0,12] = 0.5759895884974942
$arr[0,13] = 1.514644056931957
$arr[0,14] = 1.896944139870205
$arr[0,15] = 1.916393754370604
$ws.Range("Y16:AN16").Value2 = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = -1.600727426039583
$arr[0,6] = -0.7709779200558486
$arr[0,7] = -0.4516299971683568
$arr[0,8] = -0.7589144949265214
$arr[0,9] = -0.06175132635745095
$arr[0,10] = 0.5116467003986713
$arr[0,11] = 0.4136280550221194
$arr[0,12] = 0.3530477102890783
$arr[0,13] = -0.8557279162653919
$arr[0,14] = -0.701606358721174
$arr[0,15] = -0.7359525160776204
$ws.Range("AC17:AR17").Value2 = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = -0.3142152855612146
$arr[0,6] = -0.0494135395851969
$arr[0,7] = -0.2960652512497663
$arr[0,8] = 2.074687352130522
$arr[0,9] = -0.2004689067778398
$arr[0,10] = 0.2942159770784825
$arr[0,11] = 0.6923809915882817
$arr[0,12] = -0.01286797263981843
$arr[0,13] = -0.1316183744203947
$arr[0,14] = -0.1754728623905355
$arr[0,15] = -0.187152549496028
$ws.Range("AG18:AV18").Value2 = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = -0.06115297969078171
$arr[0,6] = 0.3913752358667866
$arr[0,7] = 1.196926293374756
$arr[0,8] = 0.3229362634639799
$arr[0,9] = -0.09133135081734745
$arr[0,10] = -0.0331361487157622
$arr[0,11] = -0.3126391654689975
$arr[0,12] = -0.1152140120150968
$arr[0,13] = 0.3903331526556864
$arr[0,14] = 0.5695821893874298
$arr[0,15] = 0.6150340712028246
$ws.Range("AK19:AZ19").Value2 = $arr

$arr = New-Object "object[,]" 1,12
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = -0.1051772574394816
$arr[0,6] = 0.0310212408159094
$arr[0,7] = -0.3601769225510032
$arr[0,8] = -0.544405487339128
$arr[0,9] = -0.2264357368625403
$arr[0,10] = 0.316149716722669
$arr[0,11] = 0.6473947787101642
$ws.Range("AO20:AZ20").Value2 = $arr

$arr = New-Object "object[,]" 1,8
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = -0.3725667185648485
$arr[0,6] = 0.1127973091898937
$arr[0,7] = 0.4249997313001908
$ws.Range("AS21:AZ21").Value2 = $arr

$arr = New-Object "object[,]" 1,4
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$ws.Range("AW22:AZ22").Value2 = $arr

